# Update type-declaration strings for the /Julia rule/ (see commit message:
# "chage type declation to follow julia rule").
#
# The header row of the "mergeA" sheet used leading-slash path style keys,
# one of which also encoded its element type using a non-Julia style
# ("(Int)"). Update them to the new, slash-less naming, and rewrite the
# "TEL" column type annotation using Julia's `::Vector{Int}` syntax.

$wb = $excel.ActiveWorkbook

# The author also left the selection on the "Missing" sheet at B7 before
# saving.
$wsMissing = $wb.Worksheets.Item("Missing")
$wsMissing.Activate()
$wsMissing.Range("B7").Select()

# Now edit the "mergeA" header row and leave it as the active sheet/cell.
$ws = $wb.Worksheets.Item("mergeA")
$ws.Activate()
$ws.Range("C1").Value = "Address/City"
$ws.Range("B1").Value = "Address/State"
$ws.Range("A1").Value = "Key"
$ws.Range("D1").Value = "Address/TEL::Vector{Int}"
$ws.Range("D2").Select()
